# New crime data collected — weekly refresh of the 116th Precinct CompStat sheet.
# Updates: report header (volume/week-of dates) + the weekly crime-stat grid
# (rows 15-28: Rape, Robbery, Fel. Assault, Burglary, Gr. Larceny, G.L.A.,
# TOTAL, Petit Larceny, Retail Theft, Misd. Assault, UCR Rape*, Other Sex Crimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: bump the report volume/number and the "week covering" date range.
# ---------------------------------------------------------------------------

$volCell = $ws.Range("A8")
$volText = $volCell.Text
$oldNum = "45"
$newNum = "46"
$idx = $volText.LastIndexOf($oldNum) + 1
$volCell.Characters($idx, $oldNum.Length).Text = $newNum

$weekCell = $ws.Range("C9")

$oldStart = "11/3/2025"
$newStart = "11/10/2025"
$t = $weekCell.Text
$idx = $t.IndexOf($oldStart) + 1
$weekCell.Characters($idx, $oldStart.Length).Text = $newStart

$oldEnd = "11/9/2025"
$newEnd = "11/16/2025"
$t = $weekCell.Text
$idx = $t.IndexOf($oldEnd) + 1
$weekCell.Characters($idx, $oldEnd.Length).Text = $newEnd

# ---------------------------------------------------------------------------
# Helpers for the data grid.
# ---------------------------------------------------------------------------

$fmtInt = "#,##0"
$fmtDec = "#,##0.0;""-""#,##0.0"

function Set-IntCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $c.NumberFormat = $fmtInt
}

function Set-DecCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $c.NumberFormat = $fmtDec
}

function Set-PlainCell($row, $col, $value) {
    # Plain value replacement - keeps whatever style/number-format the cell
    # already has (used where the diff shows no style ("s=...") change).
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-BlankMarker($row, $col, $text, $formatSourceRow, $formatSourceCol) {
    # Turn a numeric cell back into one of the "no data" placeholder strings
    # ("0" / "***.*") while reusing the existing placeholder style (s="13")
    # from a cell that already carries it, instead of minting a new style.
    $dst = $ws.Cells.Item($row, $col)
    $src = $ws.Cells.Item($formatSourceRow, $formatSourceCol)
    $dst.Value = "'" + $text
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-IntCell 15 4 1
Set-DecCell 15 5 -100
Set-IntCell 15 10 13
Set-DecCell 15 11 107.692307692308

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-IntCell 16 3 4
Set-IntCell 16 4 1
Set-DecCell 16 5 300
Set-IntCell 16 6 12
Set-DecCell 16 8 100
Set-IntCell 16 9 71
Set-IntCell 16 10 72
Set-DecCell 16 11 -1.388888888888
Set-DecCell 16 12 -6.578947368421

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-IntCell 17 3 9
Set-IntCell 17 4 1
Set-DecCell 17 5 800
Set-IntCell 17 6 24
Set-IntCell 17 7 14
Set-DecCell 17 8 71.428571428571
Set-IntCell 17 9 232
Set-IntCell 17 10 199
Set-DecCell 17 11 16.582914572864
Set-DecCell 17 12 14.851485148514

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-IntCell 18 3 3
Set-IntCell 18 4 2
Set-DecCell 18 5 50
Set-IntCell 18 6 7
Set-IntCell 18 7 4
Set-DecCell 18 8 75
Set-IntCell 18 9 53
Set-IntCell 18 10 63
Set-DecCell 18 11 -15.873015873015
Set-DecCell 18 12 -23.188405797101

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-IntCell 19 3 5
Set-IntCell 19 4 5
Set-IntCell 19 7 28
Set-DecCell 19 8 -14.285714285714
Set-IntCell 19 9 210
Set-IntCell 19 10 247
Set-DecCell 19 11 -14.979757085020
Set-DecCell 19 12 -21.933085501858

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-IntCell 20 3 7
Set-IntCell 20 4 5
Set-DecCell 20 5 40
Set-IntCell 20 6 16
Set-IntCell 20 7 15
Set-DecCell 20 8 6.666666666666
Set-IntCell 20 9 148
Set-IntCell 20 10 163
Set-DecCell 20 11 -9.202453987730
Set-DecCell 20 12 -18.681318681318

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold summary row; styles unchanged, values only)
# ---------------------------------------------------------------------------
Set-PlainCell 21 3 28
Set-PlainCell 21 4 15
Set-PlainCell 21 5 86.666666666666
Set-PlainCell 21 6 84
Set-PlainCell 21 7 69
Set-PlainCell 21 8 21.739130434782
Set-PlainCell 21 9 741
Set-PlainCell 21 10 761
Set-PlainCell 21 11 -2.628120893561
Set-PlainCell 21 12 -8.856088560885

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-PlainCell 24 3 11
Set-PlainCell 24 4 18
Set-PlainCell 24 5 -38.888888888888
Set-PlainCell 24 9 513
Set-PlainCell 24 10 511
Set-PlainCell 24 11 0.391389432485
Set-PlainCell 24 12 -3.389830508474

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
Set-PlainCell 25 3 3
Set-PlainCell 25 4 4
Set-PlainCell 25 5 -25
Set-PlainCell 25 6 10
Set-PlainCell 25 7 12
Set-PlainCell 25 8 -16.666666666666
Set-PlainCell 25 9 95
Set-PlainCell 25 10 139
Set-PlainCell 25 11 -31.654676258992
Set-PlainCell 25 12 -34.482758620689

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
Set-PlainCell 26 3 14
Set-PlainCell 26 4 6
Set-PlainCell 26 5 133.333333333333
Set-PlainCell 26 6 43
Set-PlainCell 26 7 29
Set-PlainCell 26 8 48.275862068965
Set-PlainCell 26 9 409
Set-PlainCell 26 10 388
Set-PlainCell 26 11 5.412371134020
Set-PlainCell 26 12 14.887640449438

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-IntCell 27 4 1
Set-DecCell 27 5 -100
Set-IntCell 27 10 22
Set-DecCell 27 11 50

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# (row 22 is an untouched "no data" row - safe to use as the format donor for
# the placeholder strings, unlike row 27 whose own D/E cells are rewritten to
# numbers earlier in this script.)
# ---------------------------------------------------------------------------
Set-IntCell 28 3 3
Set-BlankMarker 28 4 "0" 22 3
Set-BlankMarker 28 5 "***.*" 22 5
Set-IntCell 28 6 4
Set-IntCell 28 7 2
Set-DecCell 28 8 100
Set-IntCell 28 9 20
Set-DecCell 28 11 -9.090909090909
Set-DecCell 28 12 -16.666666666666
